$d = $word.ActiveDocument
$r = $d.Content
$found = $r.Find.Execute("(15)", $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)

$inner = $d.Range($r.Start + 1, $r.End - 1)
$inner.Text = "07"
